$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new RMSE values for Germany, Italy, Japan under "Univariate using Prophet" (column E)
$ws.Range("E5").Value = 0.1846
$ws.Range("E6").Value = 0.4779
$ws.Range("E7").Value = 0.4308

# Update the view: scroll so column B is the top-left visible column,
# and move the selection to E8
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("E8").Select()
